$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-5 from 45184 (2023-09-15)
# to 45185 (2023-09-16), matching the automatic daily update reflected in the diff.
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45185
}
